$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 423.64706
$ws.Range("I19").Value = 426.33334
$ws.Range("J19").Value = 422.18182
$ws.Range("K19").Value = 426.33334
$ws.Range("L19").Value = 422.18182
$ws.Range("M19").Value = -251.33334
$ws.Range("N19").Value = -772.18182
$ws.Range("H62").Value = 3684.25
$ws.Range("I62").Value = 4050
$ws.Range("J62").Value = 3318.5
$ws.Range("K62").Value = 4050
$ws.Range("L62").Value = 3318.5
$ws.Range("M62").Value = -3426
$ws.Range("N62").Value = -4566.5
$ws.Range("H65").Value = 3684.25
$ws.Range("I65").Value = 4050
$ws.Range("J65").Value = 3318.5
$ws.Range("K65").Value = 20250
$ws.Range("L65").Value = 16592.5
$ws.Range("M65").Value = -17130
$ws.Range("N65").Value = -22832.5
$ws.Range("H116").Value = 2694.0454
$ws.Range("I116").Value = 2087.3635
$ws.Range("J116").Value = 3300.7273
$ws.Range("K116").Value = 2087.3635
$ws.Range("L116").Value = 3300.7273
$ws.Range("M116").Value = 1354.6365
$ws.Range("N116").Value = -10184.7273
$ws.Range("H124").Value = 40750
$ws.Range("J124").Value = 40750
$ws.Range("L124").Value = 40750
$ws.Range("N124").Value = -50570
$ws.Range("H125").Value = 3666.6667
$ws.Range("J125").Value = 3228.5715
$ws.Range("L125").Value = 29057.1435
$ws.Range("N125").Value = -33977.1435
$ws.Range("H126").Value = 29999.6
$ws.Range("J126").Value = 29999.6
$ws.Range("L126").Value = 29999.6
$ws.Range("N126").Value = -39879.6
$ws.Range("H132").Value = 6949945
$ws.Range("I132").Value = 10421907
$ws.Range("K132").Value = 31265721
$ws.Range("M132").Value = -31263191
$ws.Range("H138").Value = 1392.5613
$ws.Range("J138").Value = 1531.026
$ws.Range("L138").Value = 4593.078
$ws.Range("N138").Value = -14873.078

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1769.2222
$ws.Range("I74").Value = 984.8333
$ws.Range("J74").Value = 3338
$ws.Range("K74").Value = 984.8333
$ws.Range("L74").Value = 3338
$ws.Range("M74").Value = -110.8333
$ws.Range("N74").Value = -5086
$ws.Range("H77").Value = 1769.2222
$ws.Range("I77").Value = 984.8333
$ws.Range("J77").Value = 3338
$ws.Range("K77").Value = 4924.1665
$ws.Range("L77").Value = 16690
$ws.Range("M77").Value = -556.1665000000003
$ws.Range("N77").Value = -25426
$ws.Range("H114").Value = 19853.908
$ws.Range("I114").Value = 8000
$ws.Range("J114").Value = 21039.3
$ws.Range("K114").Value = 8000
$ws.Range("L114").Value = 21039.3
$ws.Range("M114").Value = -3661
$ws.Range("N114").Value = -29717.3
$ws.Range("H122").Value = 1070.8438
$ws.Range("I122").Value = 929.4138
$ws.Range("K122").Value = 2788.2414
$ws.Range("M122").Value = -338.2413999999999
$ws.Range("H132").Value = 1515.04
$ws.Range("I132").Value = 1288.2285
$ws.Range("K132").Value = 3864.6855
$ws.Range("M132").Value = -1334.6855

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2636.4
$ws.Range("I86").Value = 3394.8572
$ws.Range("K86").Value = 3394.8572
$ws.Range("M86").Value = -2271.8572
$ws.Range("H89").Value = 2636.4
$ws.Range("I89").Value = 3394.8572
$ws.Range("K89").Value = 16974.286
$ws.Range("M89").Value = -11358.286
$ws.Range("H110").Value = 48999.332
$ws.Range("J110").Value = 48999.332
$ws.Range("L110").Value = 48999.332
$ws.Range("N110").Value = -57179.332
$ws.Range("H134").Value = 4128.39
$ws.Range("I134").Value = 978.3333
$ws.Range("J134").Value = 17122.375
$ws.Range("K134").Value = 2934.9999
$ws.Range("L134").Value = 51367.125
$ws.Range("M134").Value = -399.9998999999998
$ws.Range("N134").Value = -56437.125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2053.3076
$ws.Range("I31").Value = 2151.7144
$ws.Range("J31").Value = 1640
$ws.Range("K31").Value = 2151.7144
$ws.Range("L31").Value = 1640
$ws.Range("M31").Value = -1856.7144
$ws.Range("N31").Value = -2230
$ws.Range("H34").Value = 2053.3076
$ws.Range("I34").Value = 2151.7144
$ws.Range("J34").Value = 1640
$ws.Range("K34").Value = 2151.7144
$ws.Range("L34").Value = 1640
$ws.Range("M34").Value = -1949.7144
$ws.Range("N34").Value = -2044
$ws.Range("H58").Value = 669.0599999999999
$ws.Range("I58").Value = 635.1212
$ws.Range("J58").Value = 734.94116
$ws.Range("K58").Value = 635.1212
$ws.Range("L58").Value = 734.94116
$ws.Range("M58").Value = -432.1212
$ws.Range("N58").Value = -1140.94116
$ws.Range("H132").Value = 2094.0908
$ws.Range("I132").Value = 1682.3928
$ws.Range("J132").Value = 4399.6
$ws.Range("K132").Value = 5047.178400000001
$ws.Range("L132").Value = 13198.8
$ws.Range("M132").Value = -2517.178400000001
$ws.Range("N132").Value = -18258.8
$ws.Range("H134").Value = 1043.5883
$ws.Range("I134").Value = 1066.0741
$ws.Range("J134").Value = 956.8570999999999
$ws.Range("K134").Value = 3198.2223
$ws.Range("L134").Value = 2870.5713
$ws.Range("M134").Value = -663.2223000000004
$ws.Range("N134").Value = -7940.5713
$ws.Range("H136").Value = 669.0599999999999
$ws.Range("I136").Value = 635.1212
$ws.Range("J136").Value = 734.94116
$ws.Range("K136").Value = 1905.3636
$ws.Range("L136").Value = 2204.82348
$ws.Range("M136").Value = 644.6363999999999
$ws.Range("N136").Value = -7304.82348

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3787.2144
$ws.Range("I104").Value = 2990.8
$ws.Range("J104").Value = 4229.6665
$ws.Range("K104").Value = 8972.400000000001
$ws.Range("L104").Value = 12688.9995
$ws.Range("M104").Value = -6351.400000000001
$ws.Range("N104").Value = -17930.9995
$ws.Range("H113").Value = 708.86206
$ws.Range("J113").Value = 709.1786
$ws.Range("L113").Value = 2127.5358
$ws.Range("N113").Value = -6467.5358
$ws.Range("H140").Value = 22357.59
$ws.Range("J140").Value = 2928.4119
$ws.Range("L140").Value = 8785.235700000001
$ws.Range("N140").Value = -19145.2357

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7224.9473
$ws.Range("I102").Value = 5620
$ws.Range("J102").Value = 10702.333
$ws.Range("K102").Value = 5620
$ws.Range("L102").Value = 10702.333
$ws.Range("M102").Value = -3998
$ws.Range("N102").Value = -13946.333
$ws.Range("H126").Value = 1999.9
$ws.Range("I126").Value = 1749.875
$ws.Range("K126").Value = 5249.625
$ws.Range("M126").Value = -2779.625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2151
$ws.Range("I7").Value = 2101.3333
$ws.Range("K7").Value = 2101.3333
$ws.Range("M7").Value = -1989.3333
$ws.Range("H122").Value = 22729868
$ws.Range("I122").Value = 31252444
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 93757332
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -93754882
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 2151
$ws.Range("I126").Value = 2101.3333
$ws.Range("K126").Value = 6303.999899999999
$ws.Range("M126").Value = -3833.999899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 475.25
$ws.Range("I81").Value = 475.25
$ws.Range("K81").Value = 950.5
$ws.Range("M81").Value = 110.5
$ws.Range("H84").Value = 475.25
$ws.Range("I84").Value = 475.25
$ws.Range("K84").Value = 4752.5
$ws.Range("M84").Value = 551.5
